# Merge the split "<id>...</id>" runs into a single run per occurrence,
# renumbering "p042v_aN" -> "p042v_N" (dropping the stray "a").
#
# Each of the 6 occurrences in this document currently looks like three
# separate runs: "<id>", "p042v_aN", "</id>". Find & Replace across the
# run boundaries collapses them into a single run (taking on the
# formatting of the first/matched run), which is exactly what the target
# diff shows.

$d = $word.ActiveDocument

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p042v_a$i</id>"
    $new = "<id>p042v_$i</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $new, 2)
    Write-Output ("Replace #" + $i + ": " + $old + " -> " + $new + " : found=" + $found)
}
